$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")
$ws.Activate()

# Insert a new formatted row 47, inheriting formatting from row 46 above
# (mirrors dragging the bottom of the table down one row in Excel).
$ws.Rows.Item(47).Insert()

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 45233
$ws.Range("C47").Value = "2023-11-03"
$ws.Range("D47").Value = 1
$ws.Range("E47").Formula = "=G46"
$ws.Range("F47").Value = 5840
$ws.Range("G47").Formula = "=E47+F47"
$ws.Range("H47").Value = "ESPORTS"
$ws.Range("I47").Value = "EUROPEAN CIRCUIT"
$ws.Range("J47").Value = "RSA"
$ws.Range("K47").Value = "GANA MAPA 2"
$ws.Range("L47").Value = 1
$ws.Range("M47").Value = 0
$ws.Range("N47").Formula = "=ROUND((G47/`$E`$31-1)*100, 3)+`$N`$29"

# Match the scroll/selection state left behind in the saved workbook.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H43").Select()
